# Update the Files cell (I2) to reflect the new FASTQ naming convention:
# "test1.fastq.gz, test2.fastq.gz" -> "test1_R1.fastq.gz, test1_R2.fastq.gz"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "test1_R1.fastq.gz, test1_R2.fastq.gz"

# Reflect the updated active-cell selection on the (frozen) sheet view.
$ws.Range("H3").Select()
